# Fruta / hortaliza, semanal
# Insert a new weekly record for "Vega Modelo de Temuco" / Espinaca.
# The new record is inserted at row 82, pushing the previous rows
# 82..93 down to 83..94 (dimension grows from A1:R93 to A1:R94).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 82, shifting 82..93 -> 83..94
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with the new weekly entry.
$ws.Range("A82").Value = 10
$ws.Range("B82").Value = "Vega Modelo de Temuco"
$ws.Range("C82").Value = "La Araucanía"
$ws.Range("D82").Value = 44504
$ws.Range("E82").Value = 9
$ws.Range("F82").Value = 100112012
$ws.Range("G82").Value = "Espinaca"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 95
$ws.Range("K82").Value = 8000
$ws.Range("L82").Value = 8000
$ws.Range("M82").Value = 8000
$ws.Range("N82").Value = "$/docena de atados"
$ws.Range("O82").Value = "Región de La Araucanía"
$ws.Range("P82").Value = 2667
$ws.Range("Q82").Value = 3
$ws.Range("R82").Value = "Hortaliza"
